$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns A (Designator) and B (Value) for rows 1-77 using a bounded 3-way
# cut rotation through a scratch column (E) so only the used range is touched
# (avoids materializing the whole 1,048,576-row column that a full Columns.Cut
# would otherwise force).
$ws.Range("A1:A77").Cut($ws.Range("E1:E77"))
$ws.Range("B1:B77").Cut($ws.Range("A1:A77"))
$ws.Range("E1:E77").Cut($ws.Range("B1:B77"))

# Remove the now-empty scratch column E and shift it back out of the used range.
$ws.Range("E1:E77").Delete(-4159)

# The swap also moved the per-column custom widths (18.88671875 / 19) along
# with the cells, so column A/B already carry the right values - nothing else
# to do there.

# Rename the header cells for the new layout: A=Comment, D=JLCPCB Part #
$ws.Range("A1").Value = "Comment"
$ws.Range("D1").Value = "JLCPCB Part #"

Write-Output "done"
